# Add two new rows (10 and 11) to the sheet, matching the style of the
# existing data rows, with the new "日落" (Sunset) and "大海" (Ocean) entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row (row 9) onto the
# two new rows so the new cells get the same cell style ("s=1") as the rest
# of the table. Only copy A:E and H, since F/G are left blank in the source
# row and should remain entirely absent (no empty cells) in the new rows.
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Range("H9").Copy()
$ws.Range("H10").PasteSpecial(-4122)

$ws.Range("A9:E9").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("H9").Copy()
$ws.Range("H11").PasteSpecial(-4122)

# Row 10: "日落" (Sunset) entry
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "日落"
$ws.Range("C10").Value = "玩法: 上載一張有天空的相，黄昏的太陽就可出現"
$ws.Range("D10").Value = "00009.png"
$ws.Range("E10").Value = "night.jpg"
$ws.Range("H10").Value = "图轉為黄昏，見到淡橙紅色的日落，大大的太陽佔图上的天空的一半。"

# Row 11: "大海" (Ocean) entry
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "大海"
$ws.Range("C11").Value = "玩法: 上載一張自己的相, ，地就可以轉為大海"
$ws.Range("D11").Value = "00010.png"
$ws.Range("E11").Value = "boat.png"
$ws.Range("H11").Value = "將地下變成大海，相的主体不變，周邊建筑物消失，轉為海島背景，所有主体上的人物不變。大海有不同生物。相片要有真實感。"
